$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, shifting existing rows 10-35 down to 11-36
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly price record
$ws.Cells.Item(10, 1).Value = 10
$ws.Cells.Item(10, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(10, 3).Value = "La Araucanía"
$ws.Cells.Item(10, 4).Value = 44715
$ws.Cells.Item(10, 5).Value = 9
$ws.Cells.Item(10, 6).Value = 100112010
$ws.Cells.Item(10, 7).Value = "Achicoria"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 25
$ws.Cells.Item(10, 11).Value = 11000
$ws.Cells.Item(10, 12).Value = 11000
$ws.Cells.Item(10, 13).Value = 11000
$ws.Cells.Item(10, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(10, 15).Value = "Región Metropolitana"
$ws.Cells.Item(10, 16).Value = 611
$ws.Cells.Item(10, 17).Value = 18
$ws.Cells.Item(10, 18).Value = "Hortaliza"
